# Applies the "Corrección documentos individuales student 3" edit:
#   - Group code C1.010 -> C2.010 (both occurrences)
#   - Member list: remove Alejandro Soult Toscano, remove the duplicate
#     Mario Benítez Galván entry's predecessor slot (keep Mario once),
#     remove Marta de la Calle González and Manuel Alcaraz Zambrano
#   - Fecha: 20/02/2025 -> 02/07/2025
#   - Version-table dates 20/2/2025 -> 03/07/2025 (two rows)

$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Member list paragraphs: delete Manuel, Marta and Alejandro
#    (highest paragraph index first so earlier indices stay valid).
#    Mario Benítez Galván and Enrique Nicolae Barac Ploae keep their
#    existing paragraphs/hyperlinks untouched.
# ------------------------------------------------------------------
$targets = @(
    "Manuel Alcaraz Zambrano (manalczam@alum.us.es)",
    "Marta de la Calle González (mardegon7@alum.us.es)",
    "Alejandro Soult Toscano (alesoutos@alum.us.es)"
)

foreach ($target in $targets) {
    $paras = $d.Paragraphs
    for ($i = $paras.Count; $i -ge 1; $i--) {
        $p = $paras.Item($i)
        $text = $p.Range.Text
        if ($text.TrimEnd([char]13, [char]7) -eq $target) {
            $p.Range.Delete()
            break
        }
    }
}

# ------------------------------------------------------------------
# 2) Group code: C1.010 -> C2.010 (appears in the "Grupo:" line and in
#    the "Soy Mario Benítez Galván..." introduction paragraph).
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("C1.010", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "C2.010", 2) | Out-Null

# ------------------------------------------------------------------
# 3) Header date: Fecha: 20/02/2025 -> 02/07/2025
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("20/02/2025", $true, $false, $false, $false, $false, `
                   $true, 1, $false, "02/07/2025", 2) | Out-Null

# ------------------------------------------------------------------
# 4) Version table: both "20/2/2025" rows -> "03/07/2025"
# ------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $cell = $tbl.Cell($r, 1)
    $cellRng = $cell.Range
    $cellRng.Find.Execute("20/2/2025", $true, $false, $false, $false, `
                           $false, $true, 0, $false, "03/07/2025", 2) | Out-Null
}
